$d = $word.ActiveDocument

# Right-to-left script support: introduce per-column "join" paragraph
# styles (MSC_Join_A/B/C) mirroring the existing MSC_Paragraph_A/B/C
# styles, each based on the original shared "MSCJoin" style.
$joinA = $d.Styles.Add("MSC_Join_A", 1)
$joinA.BaseStyle = $d.Styles("MSCJoin")

$joinB = $d.Styles.Add("MSC_Join_B", 1)
$joinB.BaseStyle = $d.Styles("MSCJoin")

$joinC = $d.Styles.Add("MSC_Join_C", 1)
$joinC.BaseStyle = $d.Styles("MSCJoin")

# The document body only ever used the shared "MSCJoin" style for the
# first ("A") column of Scripture text, so repoint those paragraphs at
# the new column-specific style.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "MSC_Join") {
        $p.Style = $d.Styles("MSC_Join_A")
    }
}
